$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.218.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.895.46'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.893.91'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.139'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.05%  '
$ws.Range("E12").Value = '  -2.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000219'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.31%  '
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.374.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.127.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.900.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.659'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.09'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.55%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.25%  '
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.32%  '
$ws.Range("E34").Value = '  -4.46%  '
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.972'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("E39").Value = '  -6.91%  '
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("E43").Value = '  -5.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.670.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0333'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '357.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.32%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.30%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.102'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.51%  '
